$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect using the known password so we can
# write the updated holdings values, then re-protect it afterwards.
$sheetPassword = "D382"
$ws.Unprotect($sheetPassword)

# Update the "as of" date in the confidentiality / disclosure note (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Updated Weight (column D) and Percent Change (column E) values for each
# holding row (rows 2-56), per the latest model holdings snapshot.
$ws.Range("D2").Value = 0.01553455110610543
$ws.Range("E2").Value = -0.006792086679963227
$ws.Range("D3").Value = 0.05251548980619491
$ws.Range("E3").Value = -0.02203461401037643
$ws.Range("D4").Value = 0.01449565161900338
$ws.Range("E4").Value = -0.00449313192691192
$ws.Range("D5").Value = 0.009645572453913414
$ws.Range("E5").Value = 0.004758842443729927
$ws.Range("D6").Value = 0.01537355378582162
$ws.Range("E6").Value = -0.01255230125522999
$ws.Range("D7").Value = 0.01967882868812881
$ws.Range("E7").Value = -0.003407407407407526
$ws.Range("D8").Value = 0.004157086654124561
$ws.Range("E8").Value = 0.009452684353457252
$ws.Range("D9").Value = 0.006567425267720493
$ws.Range("E9").Value = 0.006290377422645488
$ws.Range("D10").Value = 0.01434070216569228
$ws.Range("E10").Value = 0.001265182186234615
$ws.Range("D11").Value = 0.008357842009262338
$ws.Range("E11").Value = -0.0001484340210776702
$ws.Range("D12").Value = 0.01506198008526852
$ws.Range("E12").Value = 0.021003212255992
$ws.Range("D13").Value = 0.003031842235473652
$ws.Range("E13").Value = -0.04475474400286428
$ws.Range("D14").Value = 0.006164792400850062
$ws.Range("E14").Value = -0.02767017155506368
$ws.Range("D15").Value = 0.01379360281493012
$ws.Range("E15").Value = 0.01382368283776714
$ws.Range("D16").Value = 0.01000828939876395
$ws.Range("E16").Value = 0.008645935016036832
$ws.Range("D17").Value = 0.02158530244614186
$ws.Range("E17").Value = -0.02838484417431775
$ws.Range("D18").Value = 0.008269015901519546
$ws.Range("E18").Value = 0.007403906742281086
$ws.Range("D19").Value = 0.01639660475998493
$ws.Range("E19").Value = 0.01549543005871312
$ws.Range("D20").Value = 0.01153250694939924
$ws.Range("E20").Value = 0.01226333907056798
$ws.Range("D21").Value = 0.007229527135071437
$ws.Range("E21").Value = 0.01650793650793658
$ws.Range("D22").Value = 0.01349338049546426
$ws.Range("E22").Value = 0.04414982623246244
$ws.Range("D23").Value = 0.01955824352510369
$ws.Range("E23").Value = 0.01036136324713133
$ws.Range("D24").Value = 0.009702112256432657
$ws.Range("E24").Value = 0.01334940205803292
$ws.Range("D25").Value = 0.02064661146257227
$ws.Range("E25").Value = -0.01382596469301689
$ws.Range("D26").Value = 0.01382077169425368
$ws.Range("E26").Value = 0.009943449575871899
$ws.Range("D27").Value = 0.02105786641671662
$ws.Range("E27").Value = -0.01774026310673327
$ws.Range("D28").Value = 0.05754964123017226
$ws.Range("E28").Value = -0.03538554398672089
$ws.Range("D29").Value = 0.020980825895895
$ws.Range("E29").Value = 0.004017857142857295
$ws.Range("D30").Value = 0.02979644491293793
$ws.Range("E30").Value = -0.01240736114580743
$ws.Range("D31").Value = 0.01577069705036298
$ws.Range("E31").Value = -0.02408498872153086
$ws.Range("D32").Value = 0.01363828118518993
$ws.Range("E32").Value = -0.01391289319045963
$ws.Range("D33").Value = 0.01898546400073432
$ws.Range("E33").Value = -0.007174782404140267
$ws.Range("D34").Value = 0.04360195737359893
$ws.Range("E34").Value = -0.01547108933540464
$ws.Range("D35").Value = 0.01084522114367988
$ws.Range("E35").Value = -0.002059025394646641
$ws.Range("D36").Value = 0.009958945007207115
$ws.Range("E36").Value = -0.006801534705266787
$ws.Range("D37").Value = 0.0108326601891981
$ws.Range("E37").Value = -0.04101352802233216
$ws.Range("D38").Value = 0.007303962420889419
$ws.Range("E38").Value = -0.009171974522292903
$ws.Range("D39").Value = 0.01169899387201213
$ws.Range("E39").Value = 0.01084812623274156
$ws.Range("D40").Value = 0.01682243662426863
$ws.Range("E40").Value = -0.008989642368575357
$ws.Range("D41").Value = 0.01701367327941599
$ws.Range("E41").Value = 0.002147409254640964
$ws.Range("D42").Value = 0.03312742395329259
$ws.Range("E42").Value = -0.02376137512639043
$ws.Range("D43").Value = 0.01127628456796849
$ws.Range("E43").Value = -0.0009260168940600222
$ws.Range("D44").Value = 0.02236442278573916
$ws.Range("E44").Value = -0.002493443961996578
$ws.Range("D45").Value = 0.01256486233428979
$ws.Range("E45").Value = -0.01244057404363075
$ws.Range("D46").Value = 0.008500974860949832
$ws.Range("E46").Value = -0.0008865539320311866
$ws.Range("D47").Value = 0.0127379243738166
$ws.Range("E47").Value = 0.0229653328658459
$ws.Range("D48").Value = 0.009910593086127852
$ws.Range("E48").Value = 0.02365863962822146
$ws.Range("D49").Value = 0.01536381516926044
$ws.Range("E49").Value = 0.01309319057193381
$ws.Range("D50").Value = 0.008467634055843859
$ws.Range("E50").Value = 0.0004981320049812776
$ws.Range("D51").Value = 0.01188576441001041
$ws.Range("E51").Value = 0.008031730292513695
$ws.Range("D52").Value = 0.008635485625363409
$ws.Range("E52").Value = -0.006321112515802696
$ws.Range("D53").Value = 0.009782657438628198
$ws.Range("E53").Value = 0.02751886373723922
$ws.Range("D54").Value = 0.1346544245152058
$ws.Range("E54").Value = 0
$ws.Range("D55").Value = 0.04390937510402719
$ws.Range("E55").Value = -0.006604226705091309
$ws.Range("E56").Value = -0.005343513155953716

# Restore sheet protection.
$ws.Protect($sheetPassword)
